# Update "countries & provincias Spain" COVID data sheet:
#  - refresh the "Datos actualizados" timestamp
#  - update case counts for several countries (columns B..H)
#  - a handful of country names (column A) swap position with their
#    neighboring row because the shared-string table was reordered
#    upstream while row-level numeric data stayed keyed by row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 'Datos actualizados a 5 de Abril de 2020 a las 14:52'
$ws.Range("B7").Value = 96471
$ws.Range("C7").Value = 379
$ws.Range("E7").Value = 68624
$ws.Range("G7").Value = 3
$ws.Range("H7").Value = 1447
$ws.Range("D37").Value = 208
$ws.Range("E37").Value = 2644
$ws.Range("B39").Value = 2385
$ws.Range("C39").Value = 206
$ws.Range("D39").Value = 488
$ws.Range("E39").Value = 1863
$ws.Range("G39").Value = 5
$ws.Range("H39").Value = 34
$ws.Range("E67").Value = 792
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = 12
$ws.Range("D94").Value = 30
$ws.Range("E94").Value = 266
$ws.Range("B98").Value = 266
$ws.Range("C98").Value = 7
$ws.Range("D98").Value = 35
$ws.Range("E98").Value = 199
$ws.Range("E142").Value = 38
$ws.Range("G142").Value = 1
$ws.Range("H142").Value = 1
$ws.Range("A156").Value = 'Haiti'
$ws.Range("C156").Value = 1
$ws.Range("D156").Value = 1
$ws.Range("H156").Value = 0
$ws.Range("A157").Value = 'Birmania'
$ws.Range("C157").Value = 0
$ws.Range("D157").Value = 0
$ws.Range("H157").Value = 1
$ws.Range("A173").Value = 'Liberia'
$ws.Range("B173").Value = 13
$ws.Range("C173").Value = 3
$ws.Range("D173").Value = 3
$ws.Range("E173").Value = 7
$ws.Range("G173").Value = 2
$ws.Range("H173").Value = 3
$ws.Range("A174").Value = 'Fiyi'
$ws.Range("F174").Value = 0
$ws.Range("A175").Value = 'Granada'
$ws.Range("B175").Value = 12
$ws.Range("C175").Value = 0
$ws.Range("E175").Value = 12
$ws.Range("F175").Value = 2
$ws.Range("A176").Value = 'Laos'
$ws.Range("C176").Value = 1
$ws.Range("D176").Value = 0
$ws.Range("E176").Value = 11
$ws.Range("A177").Value = 'Groenlandia'
$ws.Range("D177").Value = 3
$ws.Range("E177").Value = 8
$ws.Range("H177").Value = 0
$ws.Range("A178").Value = 'Curazao'
$ws.Range("B178").Value = 11
$ws.Range("D178").Value = 5
$ws.Range("E178").Value = 5
$ws.Range("H178").Value = 1
$ws.Range("A179").Value = 'Seychelles'
$ws.Range("E179").Value = 10
$ws.Range("H179").Value = 0
$ws.Range("A181").Value = 'Surinam'
$ws.Range("D181").Value = 0
$ws.Range("E181").Value = 9
$ws.Range("H181").Value = 1
$ws.Range("A183").Value = 'Angola'
$ws.Range("D183").Value = 2
$ws.Range("H183").Value = 2
$ws.Range("A184").Value = 'San Cristobal y Nieves'
$ws.Range("A186").Value = 'Republica del Chad'
$ws.Range("A193").Value = 'Somalia'
$ws.Range("D193").Value = 1
$ws.Range("H193").Value = 0
$ws.Range("A194").Value = 'Cabo Verde'
$ws.Range("D194").Value = 0
$ws.Range("H194").Value = 1
